$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new data point ("2026/01/12", 月, 7, 139) was recorded between the existing
# "2026/01/12" rows and the "2026/12/29" rows. Insert a row at 633 (shifting
# 633:674 down to 634:675) and fill it in.
$ws.Rows.Item(633).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel
# dates. Force text formatting before assigning so Excel doesn't
# auto-convert the "yyyy/mm/dd" string into a date serial number, then
# clear the format again so the cell's style matches its untouched
# neighbours (format cleanup does not revert the already-committed text
# value back to a number).
$ws.Cells.Item(633, 1).NumberFormat = "@"
$ws.Cells.Item(633, 1).Value = "2026/01/12"
$ws.Cells.Item(633, 1).ClearFormats()

$ws.Cells.Item(633, 2).Value = "月"
$ws.Cells.Item(633, 3).Value = 7
$ws.Cells.Item(633, 4).Value = 139
